$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Update status ("Estatus", column F) for tasks that are now done
$ws.Range("F9").Value = "Hecho"
$ws.Range("F10").Value = "Hecho"
$ws.Range("F11").Value = "Hecho"
$ws.Range("F12").Value = "Hecho"

# Set estimated hours (column G) for rows that were missing it
$ws.Range("G10").Value = 1
$ws.Range("G11").Value = 1
$ws.Range("G12").Value = 1

# Register the consumed hour for the day the task was completed (column W)
$ws.Range("W10").Value = 1
$ws.Range("W12").Value = 1

$wb.Save()
